$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set values in shared-string insertion order: fm38, desc1, fm39, desc2, "in progress"
$ws.Range("B43").Value = "fm38"
$ws.Range("C43").Value = "BugToFix: 3 CAT XL layers net loss error at level 3"
$ws.Range("B44").Value = "fm39"
$ws.Range("C44").Value = "BugToFix: fm24 inuring level 3 outputs wrong net loss when losses carried to level 3"
$ws.Range("H43").Value = "in progress"
$ws.Range("I43").Value = "in progress"
$ws.Range("H44").Value = "in progress"
$ws.Range("I44").Value = "in progress"

# Copy formatting from the row above (row 42) onto the new rows, preserving values
$ws.Range("B42:C42").Copy()
$ws.Range("B43:C43").PasteSpecial(-4122)
$ws.Range("B42:C42").Copy()
$ws.Range("B44:C44").PasteSpecial(-4122)

$ws.Range("H42:I42").Copy()
$ws.Range("H43:I43").PasteSpecial(-4122)
$ws.Range("H42:I42").Copy()
$ws.Range("H44:I44").PasteSpecial(-4122)

# Update the sheet view: scroll to top, select B1 (instead of previous C42 / topLeftCell A24)
$ws.Range("B1").Select()
